# Applies the "searchStudeny.py" data-refresh edit:
#  - adds an email + monthly count for Jaspreet Kaur (row 6)
#  - fixes the spelling of "Randell Holland" -> "Rondell Holland" and adds
#    his email + monthly count (row 9)
#  - bumps garytsai's (row 13) monthly count from 6 to 7
#  - appends two brand-new students (Sujay Bhaskar kashyap, jimmy Barreto)
#    as rows 15 and 16
#  - extends the Monthly_STAT SUM formula to cover the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")
$stat = $wb.Worksheets.Item("Monthly_STAT")

# Row 6 - Jaspreet Kaur: add email + count
$ws.Range("C6").Value = "jaspreet.kaur6@jjay.cuny.edu"
$ws.Range("D6").Value = 3

# Row 9 - Randell Holland -> Rondell Holland, add email + count
$ws.Range("B9").Value = "Rondell Holland"
$ws.Range("C9").Value = "rondell.holland@jjay.cuny.edu"
$ws.Range("D9").Value = 1

# Row 13 - garytsai: count 6 -> 7
$ws.Range("D13").Value = 7

# Row 15 - new student: Sujay Bhaskar kashyap
$ws.Range("A15").Value = "8FD8A841E0"
$ws.Range("B15").Value = "Sujay Bhaskar kashyap"
$ws.Range("C15").Value = "sujay.bhaskarkashyap@jjay.cuny.edu"
$ws.Range("D15").Value = 2

# Row 16 - new student: jimmy Barreto
$ws.Range("A16").Value = "8FD8AC0500"
$ws.Range("B16").Value = "jimmy Barreto"
$ws.Range("C16").Value = "j.barreto1823@yahoo.com"
$ws.Range("D16").Value = 1

# Extend the monthly total formula to include the new rows
$stat.Range("B2").Formula = "=SUM(Sheet!D2:D16)"
